$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: "proc_filepath" column -> "project_filepath" ---
$ws.Range("F1").Value = "project_filepath"
$ws.Range("F5").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Re-do per-row proc/project filepaths for EasterBush / CO2_H2O (rows 2-4) ---
# Set values first, then copy the quote-prefixed / left-aligned text formatting that
# the rest of column F already uses (F5:F11) so style 1 gets reused rather than a
# brand-new style created.
$ws.Range("F2").Value = "N:/0Peter/curr/ECsystem/eddypro/EB_1995.eddypro"
$ws.Range("F3").Value = "N:/0Peter/curr/ECsystem/eddypro/EB_2007.eddypro"
$ws.Range("F4").Value = "N:/0Peter/curr/ECsystem/eddypro/EB_2013.eddypro"
$ws.Range("F5").Copy()
$ws.Range("F2:F4").PasteSpecial(-4122)

# --- New site: EB_jasmin / EasterBush_jasmin, duplicating the EasterBush block (rows 12-15) ---

# Row 12 (CO2_H2O, 1995 interval)
$ws.Range("A12").Value = "EB_jasmin"
$ws.Range("B12").Value = "EasterBush_jasmin"
$ws.Range("C12").Value = "CO2_H2O"
$ws.Range("D12").Value = 34700
$ws.Range("E12").Value = 39082.999305555553
$ws.Range("F12").Value = "/group_workspaces/jasmin2/eddystore/stations/EasterBush/proc/processing1995.eddypro"
$ws.Range("G12").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B12 & "/raw_files"'
$ws.Range("H12").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B12 & "/output"'

# Row 13 (CO2_H2O, 2007 interval)
$ws.Range("A13").Value = "EB_jasmin"
$ws.Range("B13").Value = "EasterBush_jasmin"
$ws.Range("C13").Value = "CO2_H2O"
$ws.Range("D13").Value = 39083
$ws.Range("E13").Value = 41274.999305555553
$ws.Range("F13").Value = "/group_workspaces/jasmin2/eddystore/stations/EasterBush/proc/processing2007.eddypro"
$ws.Range("G13").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B13 & "/raw_files"'
$ws.Range("H13").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B13 & "/output"'

# Row 14 (CO2_H2O, 2013 interval)
$ws.Range("A14").Value = "EB_jasmin"
$ws.Range("B14").Value = "EasterBush_jasmin"
$ws.Range("C14").Value = "CO2_H2O"
$ws.Range("D14").Value = 41275
$ws.Range("E14").Value = 43831.5
$ws.Range("F14").Value = "/group_workspaces/jasmin2/eddystore/stations/EasterBush/proc/processing2013.eddypro"
$ws.Range("G14").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B14 & "/raw_files"'
$ws.Range("H14").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B14 & "/output"'

# Row 15 (N2O_H2O, advanced processing interval)
$ws.Range("A15").Value = "EB_jasmin"
$ws.Range("B15").Value = "EasterBush_jasmin"
$ws.Range("C15").Value = "N2O_H2O"
$ws.Range("D15").Value = 42370
$ws.Range("E15").Value = 43831.5
$ws.Range("F15").Value = "/group_workspaces/jasmin2/eddystore/stations/EasterBush/proc/processing_2018-02-26T132538_adv.eddypro"
$ws.Range("G15").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B15 & "/raw_files"'
$ws.Range("H15").Formula = '="/group_workspaces/jasmin2/eddystore/stations/" & $B15 & "/output"'

# Copy formatting (text quote-prefix / left align on A,B,F ; date format on D,E) from the
# EasterBush template rows onto the new EB_jasmin rows, after the values/formulas are set.
$ws.Range("C9").Copy()
$ws.Range("A12:B15").PasteSpecial(-4122)

$ws.Range("D2:E4").Copy()
$ws.Range("D12:E14").PasteSpecial(-4122)

$ws.Range("D5:E5").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)

$ws.Range("F5").Copy()
$ws.Range("F12:F15").PasteSpecial(-4122)

# --- Row 16: a stray quote-prefixed / left-aligned blank cell in A16 ---
$ws.Range("F5").Copy()
$ws.Range("A16").PasteSpecial(-4122)

# --- Update the active selection to F5 (single cell) like in the edited workbook ---
$null = $ws.Range("F5").Select()
